# Add a new "Swiss" worksheet with Switzerland test data, based on the
# existing "Belgium" sheet (same layout/styling), populate its market name
# and ticket id, make it the active/selected sheet, and restore a
# "select-all" selection state on the Belgium sheet (as it is no longer
# the active tab).

$wb = $excel.ActiveWorkbook

# Duplicate the Belgium sheet and move the copy to the very end of the
# workbook (after the current last sheet, i.e. Czech).
$belgium = $wb.Worksheets.Item("Belgium")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$belgium.Copy($null, $lastSheet)

# The copy becomes the new last sheet; rename it and fill in the
# Switzerland-specific values.
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2348"

# Selection on the new Swiss sheet.
$swiss.Range("B6").Select()

# Belgium is no longer the active tab; update its stored selection to a
# "select all cells" state.
$belgium.Activate()
$belgium.Cells.Select()

# Re-activate Swiss so it is the tab shown/selected when the workbook is
# reopened.
$swiss.Activate()
$swiss.Range("B6").Select()
